{"js": "// Paragraph \"Proceso de inicializaci\u00f3n de archivos local al repositorio de\n// la nube \" gets rewritten as \"1.-Proceso de inicializaci\u00f3n de archivos\n// local al repositorio de la nube.\" (split across three runs: \"1.-\", the\n// sentence, and the closing \".\"), and a brand-new paragraph \"2.- Proceso de\n// commit al repositorio. \" is added right after it -- with \"commit\" wrapped\n// in spell-check <w:proofErr/> markers, since it is flagged as a foreign\n// word by the (Spanish) proofing language.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text,items/uniqueLocalId\");\nawait context.sync();\n\nconst targetText =\n  \"Proceso de inicializaci\u00f3n de archivos local al repositorio de la nube\";\nconst target = paragraphs.items.filter(\n  (p) => p.text.trim() === targetText\n)[0];\nconst targetId = target.uniqueLocalId;\n\n// Inject the two replacement paragraphs (with their exact run/proofErr\n// layout) as raw OOXML right after the target paragraph; this preserves the\n// run boundaries and the <w:proofErr/> spell-check markers verbatim,\n// instead of letting same-formatted runs coalesce into a single run the\n// way the higher-level insertText()/insertParagraph() APIs would.\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  \"<w:p><w:r><w:t>1.-</w:t></w:r><w:r><w:t>Proceso de inicializaci\u00f3n de archivos local al repositorio de la nube</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>\" +\n  '<w:p><w:r><w:t xml:space=\"preserve\">2.- Proceso de </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>commit</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> al repositorio.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r></w:p>' +\n  \"</w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\ntarget.getRange().insertOoxml(ooxml, Word.InsertLocation.after);\nawait context.sync();\n\n// The original paragraph got pushed after the two freshly-inserted ones;\n// find it again by its (preserved) paragraph id and remove the now-\n// duplicated text, leaving only the rewritten version in place.\nparagraphs.load(\"items/uniqueLocalId\");\nawait context.sync();\nconst stale = paragraphs.items.filter((p) => p.uniqueLocalId === targetId)[0];\nstale.delete();\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Paragraph 2 currently reads:\n#   \"Proceso de inicializaci\u00f3n de archivos local al repositorio de la nube \"\n# It gets rewritten (as \"1.- <text>.\" split across three runs) and a new\n# paragraph 3 (\"2.- Proceso de commit al repositorio. \") is appended right\n# after it, with \"commit\" wrapped in spell-check proofErr markers since it\n# is flagged as a foreign word.\n$p2 = $d.Paragraphs.Item(2)\n$r = $p2.Range\n\n$xml = @'\n<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t>1.-</w:t></w:r><w:r><w:t>Proceso de inicializaci\u00f3n de archivos local al repositorio de la nube</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t xml:space=\"preserve\">2.- Proceso de </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>commit</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> al repositorio.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r></w:p>\n'@\n\n$r.InsertXML($xml)\n"}
